{"js": "// Update the \"NATIONAL\" -> \"NATINIMP\" nationality-variable figures in the\n// TSDT comparison table. Only the numeric values change; everything else\n// (labels, formatting, cell borders, etc.) stays the same.\n//\n// Strategy: locate the document's single table, then for each cell that\n// needs to change, search within that specific cell's body for the old\n// value and replace it in place with Range.insertText(..., \"Replace\").\n// Scoping the search to the owning cell (rather than the whole body)\n// avoids any ambiguity from duplicate values (e.g. \"0.06\" appears twice,\n// \"1.00\" appears twice) and keeps the edit surgical -- it only touches the\n// text run, leaving paragraph/run formatting untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, oldText, newText]\nconst edits = [\n  [1, 1, \"6,560\", \"6,825\"],\n  [1, 2, \"0.83\", \"0.82\"],\n  [2, 1, \"676\", \"734\"],\n  [3, 1, \"509\", \"549\"],\n  [3, 2, \"0.06\", \"0.07\"],\n  [4, 1, \"168\", \"187\"],\n  [5, 1, \"7,913\", \"8,295\"],\n];\n\nfor (const [rowIndex, colIndex, oldText, newText] of edits) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\n      `Could not find \"${oldText}\" in table cell (${rowIndex}, ${colIndex})`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the \"NATIONAL\" -> \"NATINIMP\" nationality-variable figures in the\n# TSDT comparison table. Only the numeric values change; everything else\n# (labels, formatting, cell borders, etc.) stays the same.\n#\n# Strategy: walk the document's single table and set each target cell's\n# Range.Text directly to the new value. Addressing cells by (row, column)\n# sidesteps any ambiguity from duplicate values elsewhere in the table\n# (e.g. \"0.06\" and \"1.00\" both appear more than once), and assigning\n# Range.Text in place (rather than rebuilding the run) preserves the\n# existing xml:space=\"preserve\" formatting on the text run.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Each entry: 1-based row, 1-based column, expected old text, new text.\n$edits = @(\n    @(2, 2, \"6,560\", \"6,825\"),\n    @(2, 3, \"0.83\",  \"0.82\"),\n    @(3, 2, \"676\",   \"734\"),\n    @(4, 2, \"509\",   \"549\"),\n    @(4, 3, \"0.06\",  \"0.07\"),\n    @(5, 2, \"168\",   \"187\"),\n    @(6, 2, \"7,913\", \"8,295\")\n)\n\nforeach ($edit in $edits) {\n    $rowIndex = $edit[0]\n    $colIndex = $edit[1]\n    $oldText  = $edit[2]\n    $newText  = $edit[3]\n\n    $cell = $tbl.Cell($rowIndex, $colIndex)\n    # Cell.Range.Text includes the trailing cell-mark characters, so trim\n    # them off before comparing against the expected plain value.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -ne $oldText) {\n        Write-Output \"WARNING: cell ($rowIndex, $colIndex) was '$current', expected '$oldText'\"\n    }\n\n    $cell.Range.Text = $newText\n}\n"}
